$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.046857595443726
$ws.Range("B1").Value = 6.31736421585083
$ws.Range("C1").Value = 6.794921398162842
$ws.Range("D1").Value = 7.268158435821533
$ws.Range("E1").Value = 4.872618198394775
